# Scheduled-runner price/profit refresh across the Mateus_Profits workbook.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H-N) for the leves whose market data moved, sheet by sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 820.93335
$ws.Range("I98").Value = 808.1429000000001
$ws.Range("K98").Value = 808.1429000000001
$ws.Range("M98").Value = 689.8570999999999

$ws.Range("H116").Value = 4633.3335

$ws.Range("H118").Value = 352.75
$ws.Range("I118").Value = 339.36365
$ws.Range("K118").Value = 1018.09095
$ws.Range("M118").Value = 638.90905

$ws.Range("H122").Value = 820.93335
$ws.Range("I122").Value = 808.1429000000001
$ws.Range("K122").Value = 2424.4287
$ws.Range("M122").Value = 25.57129999999961

$ws.Range("H141").Value = 4048.5715
$ws.Range("I141").Value = 3667.6924
$ws.Range("K141").Value = 11003.0772
$ws.Range("M141").Value = -5823.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4390.709
$ws.Range("I32").Value = 3629.4473
$ws.Range("J32").Value = 23676
$ws.Range("K32").Value = 3629.4473
$ws.Range("L32").Value = 23676
$ws.Range("M32").Value = -3342.4473
$ws.Range("N32").Value = -24250

$ws.Range("H55").Value = 24433.334
$ws.Range("J55").Value = 31650
$ws.Range("L55").Value = 31650
$ws.Range("N55").Value = -32280

$ws.Range("H122").Value = 2248.3333
$ws.Range("I122").Value = 2248.3333
$ws.Range("K122").Value = 6744.999899999999
$ws.Range("M122").Value = -4294.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 57054.4
$ws.Range("J35").Value = 57054.4
$ws.Range("L35").Value = 57054.4
$ws.Range("N35").Value = -57674.4

$ws.Range("H82").Value = 23925
$ws.Range("J82").Value = 41850
$ws.Range("L82").Value = 41850
$ws.Range("N82").Value = -42616

$ws.Range("H85").Value = 23925
$ws.Range("J85").Value = 41850
$ws.Range("L85").Value = 41850
$ws.Range("N85").Value = -44502

$ws.Range("H94").Value = 865.0714
$ws.Range("I94").Value = 743.4167
$ws.Range("J94").Value = 1595
$ws.Range("K94").Value = 743.4167
$ws.Range("L94").Value = 1595
$ws.Range("M94").Value = -292.4167
$ws.Range("N94").Value = -2497

$ws.Range("H99").Value = 1796.0625
$ws.Range("I99").Value = 1984.0454
$ws.Range("K99").Value = 1984.0454
$ws.Range("M99").Value = -486.0454

$ws.Range("H107").Value = 1388.5454
$ws.Range("I107").Value = 1086.0526
$ws.Range("J107").Value = 3304.3333
$ws.Range("K107").Value = 1086.0526
$ws.Range("L107").Value = 3304.3333
$ws.Range("M107").Value = 833.9474
$ws.Range("N107").Value = -7144.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 22406.25
$ws.Range("J41").Value = 20407.143
$ws.Range("L41").Value = 20407.143
$ws.Range("N41").Value = -21263.143

$ws.Range("H51").Value = 27500
$ws.Range("J51").Value = 27500
$ws.Range("L51").Value = 27500
$ws.Range("N51").Value = -28972

$ws.Range("H59").Value = 40000
$ws.Range("I59").Value = 40000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 40000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -38855
$ws.Range("N59").ClearContents()

$ws.Range("H60").Value = 19805.555
$ws.Range("I60").Value = 17500
$ws.Range("J60").Value = 22687.5
$ws.Range("K60").Value = 17500
$ws.Range("L60").Value = 22687.5
$ws.Range("M60").Value = -16989
$ws.Range("N60").Value = -23709.5

$ws.Range("H61").Value = 27500
$ws.Range("J61").Value = 27500
$ws.Range("L61").Value = 27500
$ws.Range("N61").Value = -28196

$ws.Range("H68").Value = 37300
$ws.Range("J68").Value = 35925
$ws.Range("L68").Value = 35925
$ws.Range("N68").Value = -37423

$ws.Range("H71").Value = 37300
$ws.Range("J71").Value = 35925
$ws.Range("L71").Value = 107775
$ws.Range("N71").Value = -115263

$ws.Range("H74").Value = 40279.75
$ws.Range("J74").Value = 40279.75
$ws.Range("L74").Value = 40279.75
$ws.Range("N74").Value = -42027.75

$ws.Range("H77").Value = 40279.75
$ws.Range("J77").Value = 40279.75
$ws.Range("L77").Value = 120839.25
$ws.Range("N77").Value = -129575.25

$ws.Range("H116").Value = 55999.668
$ws.Range("J116").Value = 55999.668
$ws.Range("L116").Value = 55999.668
$ws.Range("N116").Value = -65177.668

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H120").Value = 24754.666
$ws.Range("J120").Value = 24754.666
$ws.Range("L120").Value = 24754.666
$ws.Range("N120").Value = -32012.666

$ws.Range("H121").Value = 40326
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 40326
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 40326
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -42946

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2588.9
$ws.Range("I122").Value = 3098
$ws.Range("J122").Value = 2079.8
$ws.Range("K122").Value = 9294
$ws.Range("L122").Value = 6239.400000000001
$ws.Range("M122").Value = -6844
$ws.Range("N122").Value = -11139.4

$ws.Range("H126").Value = 5953.4165
$ws.Range("I126").Value = 3420.4285
$ws.Range("K126").Value = 10261.2855
$ws.Range("M126").Value = -7791.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3267.85
$ws.Range("I7").Value = 3297.6667
$ws.Range("K7").Value = 3297.6667
$ws.Range("M7").Value = -3185.6667

$ws.Range("H40").Value = 3004.125
$ws.Range("I40").Value = 1961.8096
$ws.Range("J40").Value = 4994
$ws.Range("K40").Value = 1961.8096
$ws.Range("L40").Value = 4994
$ws.Range("M40").Value = -1825.8096
$ws.Range("N40").Value = -5266

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H126").Value = 3267.85
$ws.Range("I126").Value = 3297.6667
$ws.Range("K126").Value = 9893.000100000001
$ws.Range("M126").Value = -7423.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 30068
$ws.Range("J51").Value = 30068
$ws.Range("L51").Value = 30068
$ws.Range("N51").Value = -31088

$ws.Range("H52").Value = 27500
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H54").Value = 79550
$ws.Range("J54").Value = 99233.336
$ws.Range("L54").Value = 99233.336
$ws.Range("N54").Value = -100273.336

$ws.Range("H70").Value = 27852.5
$ws.Range("I70").Value = 25600
$ws.Range("J70").Value = 30105
$ws.Range("K70").Value = 25600
$ws.Range("L70").Value = 30105
$ws.Range("M70").Value = -25285
$ws.Range("N70").Value = -30735

$ws.Range("H73").Value = 27852.5
$ws.Range("I73").Value = 25600
$ws.Range("J73").Value = 30105
$ws.Range("K73").Value = 25600
$ws.Range("L73").Value = 30105
$ws.Range("M73").Value = -24508
$ws.Range("N73").Value = -32289

$ws.Range("H107").Value = 1043.1428
$ws.Range("J107").Value = 788
$ws.Range("L107").Value = 2364
$ws.Range("N107").Value = -6204

$ws.Range("H122").Value = 2766.0334
$ws.Range("I122").Value = 1707.875
$ws.Range("J122").Value = 6998.6665
$ws.Range("K122").Value = 5123.625
$ws.Range("L122").Value = 20995.9995
$ws.Range("M122").Value = -2673.625
$ws.Range("N122").Value = -25895.9995

$ws.Range("H126").Value = 9158.9
$ws.Range("I126").Value = 7084.143
$ws.Range("K126").Value = 21252.429
$ws.Range("M126").Value = -18782.429
